$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 448, shifting existing rows 448:551 down to 449:552
$ws.Rows.Item(448).Insert()

# Populate the new row 448 with the new price-record data
$ws.Range("A448").Value = 11
$ws.Range("B448").Value = "Vega Monumental Concepción"
$ws.Range("C448").Value = "Bíobío"
$ws.Range("D448").Value = 45204
$ws.Range("E448").Value = 8
$ws.Range("F448").Value = "Fruta"
$ws.Range("G448").Value = 100102
$ws.Range("H448").Value = "Cítricos"
$ws.Range("I448").Value = 100102005
$ws.Range("J448").Value = "Naranja"
$ws.Range("K448").Value = "Navel Late"
$ws.Range("L448").Value = "Primera"
$ws.Range("M448").Value = 270
$ws.Range("N448").Value = 9000
$ws.Range("O448").Value = 10000
$ws.Range("P448").Value = 9444
$ws.Range("Q448").Value = "$/bandeja 15 kilos granel"
$ws.Range("R448").Value = "Región de O'Higgins"
$ws.Range("S448").Value = 630
$ws.Range("T448").Value = 15
